$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5.0
$ws.Range("B3").Value = 5.0
$ws.Range("B4").Value = 5.0
$ws.Range("B5").Value = 5.0
$ws.Range("B6").Value = 5.0
$ws.Range("B9").Value = 5.0
$ws.Range("B10").Value = 5.0
$ws.Range("B11").Value = 5.0
$ws.Range("B12").Value = 5.0
$ws.Range("B13").Value = 5.0
$ws.Range("B14").Value = 5.0
$ws.Range("B15").Value = 5.0
$ws.Range("B16").Value = 1.0
$ws.Range("B17").Value = 5.0
$ws.Range("B18").Value = 5.0
$ws.Range("B19").Value = 5.0
$ws.Range("B20").Value = 5.0
$ws.Range("B21").Value = 5.0
$ws.Range("B22").Value = 5.0
$ws.Range("B23").Value = 5.0
$ws.Range("B24").Value = 5.0
$ws.Range("B25").Value = 5.0
$ws.Range("B26").Value = 5.0
$ws.Range("B27").Value = 5.0
$ws.Range("B28").Value = 5.0
$ws.Range("B29").Value = 5.0
$ws.Range("B30").Value = 5.0
$ws.Range("B31").Value = 5.0
$ws.Range("B32").Value = 5.0
$ws.Range("B33").Value = 5.0
$ws.Range("B34").Value = 5.0
$ws.Range("B35").Value = 5.0
$ws.Range("B36").Value = 5.0
$ws.Range("B37").Value = 5.0
$ws.Range("B39").Value = 5.0
$ws.Range("B40").Value = 5.0
$ws.Range("B41").Value = 5.0
$ws.Range("B42").Value = 5.0
$ws.Range("B43").Value = 5.0
$ws.Range("B44").Value = 5.0
$ws.Range("B45").Value = 5.0
$ws.Range("B46").Value = 5.0
$ws.Range("B47").Value = 5.0
$ws.Range("B48").Value = 5.0
$ws.Range("B49").Value = 5.0
$ws.Range("B50").Value = 5.0
$ws.Range("B51").Value = 5.0
$ws.Range("B52").Value = 5.0
$ws.Range("B53").Value = 5.0
$ws.Range("B54").Value = 5.0
$ws.Range("B55").Value = 5.0
$ws.Range("B56").Value = 5.0
$ws.Range("B57").Value = 5.0
$ws.Range("B58").Value = 5.0
$ws.Range("B59").Value = 5.0
$ws.Range("B60").Value = 5.0
$ws.Range("B61").Value = 5.0
$ws.Range("B62").Value = 4.0
$ws.Range("B63").Value = 5.0
$ws.Range("B64").Value = 5.0
$ws.Range("B65").Value = 5.0
$ws.Range("B66").Value = 5.0
$ws.Range("B67").Value = 5.0
$ws.Range("B68").Value = 5.0
$ws.Range("B69").Value = 5.0
$ws.Range("B70").Value = 5.0
$ws.Range("B71").Value = 5.0
$ws.Range("B72").Value = 5.0
$ws.Range("B73").Value = 5.0
$ws.Range("B74").Value = 5.0
$ws.Range("B75").Value = 5.0
$ws.Range("B76").Value = 5.0
$ws.Range("B77").Value = 5.0
$ws.Range("B78").Value = 5.0
$ws.Range("B79").Value = 5.0
$ws.Range("B80").Value = 5.0
$ws.Range("B81").Value = 5.0
$ws.Range("B82").Value = 5.0
$ws.Range("B83").Value = 5.0
$ws.Range("B84").Value = 5.0
$ws.Range("B85").Value = 5.0
$ws.Range("B86").Value = 5.0
$ws.Range("B87").Value = 5.0
$ws.Range("B88").Value = 5.0
$ws.Range("B89").Value = 5.0
$ws.Range("B90").Value = 5.0
$ws.Range("B91").Value = 5.0
$ws.Range("B92").Value = 5.0
$ws.Range("B93").Value = 5.0
$ws.Range("B94").Value = 5.0
$ws.Range("B95").Value = 5.0
$ws.Range("B96").Value = 5.0
$ws.Range("B97").Value = 5.0
$ws.Range("B98").Value = 5.0
$ws.Range("B99").Value = 5.0
$ws.Range("B100").Value = 5.0
$ws.Range("B101").Value = 5.0
